$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Spelling / naming corrections -------------------------------------
$ws.Range("B6").Value = "Estudio e Implementación UARGFlow"

$ws.Range("Q2").Value = "Casos de Pruebas Integración Testify"
$ws.Range("Q3").Value = "Casos de Pruebas Integración Testify"
$ws.Range("Q4").Value = "Casos de Pruebas Integración Testify"

# --- New entries added under Construcción Iteración 1 -------------------
# J14 already exists (empty, but pre-formatted) so assigning the value keeps
# its existing style. I15/J15 are brand new cells on that row, so first copy
# the number/cell formatting from a same-styled neighbour before writing the
# value into them.
$ws.Range("J14").Value = "CU 18 CRUD Categoria.docx"

$ws.Range("I14").Copy() | Out-Null
$ws.Range("I15").PasteSpecial(-4122) | Out-Null
$ws.Range("I15").Value = "C114"

$ws.Range("K14").Copy() | Out-Null
$ws.Range("J15").PasteSpecial(-4122) | Out-Null
$ws.Range("J15").Value = "CU 18 CRUD Categoria.pdf"

# --- Selection / view state ---------------------------------------------
$ws.Range("B6").Select() | Out-Null
